# -----------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1. "总计" (summary) sheet: shift the existing quarterly rows down by
#    one and insert the freshly reported "2022-Q4" figures at the top
#    of the table.
# 2. Add a brand-new "2022-Q4" worksheet (positioned right after "总计",
#    i.e. immediately before "2022-Q1") holding the per-fund holdings
#    detail for that quarter.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===========================================================================
# 1. Update "总计"
# ===========================================================================
$total = $wb.Worksheets.Item("总计")

# Snapshot the current (pre-edit) quarterly rows (B:D, rows 2..6) before
# overwriting anything, so they can be shifted down by one row.
$oldB = @()
$oldC = @()
$oldD = @()
for ($r = 2; $r -le 6; $r++) {
    $oldB += $total.Cells.Item($r, 2).Value2
    $oldC += $total.Cells.Item($r, 3).Value2
    $oldD += $total.Cells.Item($r, 4).Value2
}

# Row 2 becomes the new "2022-Q4" entry.
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 13
$total.Cells.Item(2, 4).Value = 1.83

# Rows 3..6 take over what used to live in rows 2..5.
for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 3
    $total.Cells.Item($r, 2).Value = $oldB[$i]
    $total.Cells.Item($r, 3).Value = $oldC[$i]
    $total.Cells.Item($r, 4).Value = $oldD[$i]
}

# New row 7 takes over what used to be row 6 ("2021-Q1"). Column A needs a
# freshly styled index cell (copy the format from the row above so it
# matches the rest of the index column).
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)
$total.Cells.Item(7, 1).Value = 5
$total.Cells.Item(7, 2).Value = $oldB[4]
$total.Cells.Item(7, 3).Value = $oldC[4]
$total.Cells.Item(7, 4).Value = $oldD[4]

# ===========================================================================
# 2. Build the new "2022-Q4" sheet
# ===========================================================================
# Duplicate "2021-Q3" (it already has the right column layout/styles and,
# at 8 rows, is the closest existing sheet to the 14 rows we need) right
# next to itself, then move/rename/resize/refill it in place. Duplicating
# (rather than Worksheets.Add()) is what lets the new sheet inherit the
# workbook's existing cell styles (bordered/bold header + index column)
# instead of Excel's blank-sheet defaults.
$template = $wb.Worksheets.Item("2021-Q3")
$template.Copy($null, $template)
$q4 = $wb.Worksheets.Item("2021-Q3 (2)")

# Grow from 8 to 14 rows by replicating the last row's formatting.
$lastRow = $q4.Range("A8:H8")
for ($r = 9; $r -le 14; $r++) {
    $destRow = $q4.Range("A" + $r + ":H" + $r)
    $lastRow.Copy()
    $destRow.PasteSpecial(-4122)
}

# Move the sheet into position (right after "总计" / right before
# "2022-Q1"). Moving changes tab ordering, which can stale out any
# worksheet handle obtained beforehand, so re-fetch a fresh handle by
# name immediately afterwards before renaming/filling it in.
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$q4.Move($q1Sheet)
$q4 = $wb.Worksheets.Item("2021-Q3 (2)")
$q4.Name = "2022-Q4"

# Header row (B..H): 基金代码/基金名称/基金规模/股票总仓位/仓位占比/
# 持有市值(亿元)/仓位排名.
$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

$rows = @(
    @("015203", "汇添富全球移动互联灵活配置混合（QDII）D", "11.52", "92.14", "3.45", "0.3974", 7),
    @("001668", "汇添富全球移动互联灵活配置混合（QDII）A", "11.48", "92.14", "3.45", "0.3961", 7),
    @("012868", "易方达标普信息科技指数（QDII-LOF）人民币 C", "5.09", "91.36", "3.75", "0.1909", 4),
    @("161128", "易方达标普信息科技指数（QDII-LOF）人民币", "5.09", "91.36", "3.75", "0.1909", 4),
    @("000043", "嘉实美国成长股票（QDII）人民币", "12.69", "92.23", "1.47", "0.1865", 9),
    @("000044", "嘉实美国成长股票（QDII）美元现汇", "12.69", "92.23", "1.47", "0.1865", 9),
    @("003721", "易方达标普信息科技指数（QDII-LOF）美元A", "4.93", "91.36", "3.75", "0.1849", 4),
    @("006792", "鹏华香港美国互联网股票（LOF）美元现汇", "1.29", "88.46", "3.37", "0.0435", 8),
    @("160644", "鹏华香港美国互联网股票（LOF）人民币", "1.29", "88.46", "3.37", "0.0435", 8),
    @("012869", "易方达标普信息科技指数（QDII-LOF）美元 C", "0.16", "91.36", "3.75", "0.0060", 4),
    @("014002", "浦银安盛全球智能科技股票（QDII）C", "0.30", "42.55", "1.25", "0.0038", 8),
    @("006555", "浦银安盛全球智能科技股票（QDII）A", "0.25", "42.55", "1.25", "0.0031", 8),
    @("015202", "汇添富全球移动互联灵活配置混合（QDII）C", "0.01", "92.14", "3.45", "0.0003", 7)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    # Column A: numeric row index (0-based), using the style already on
    # this cell from the template.
    $q4.Cells.Item($r, 1).Value = $i

    # Column B: fund code. Leading apostrophe forces text (so leading
    # zeros survive); ".Style = 'Normal'" drops the resulting
    # quote-prefix flag so the cell ends up with the default style,
    # exactly like the rest of the (script-authored) workbook.
    $q4.Cells.Item($r, 2).Value = "'" + $data[0]
    $q4.Cells.Item($r, 2).Style = "Normal"

    # Column C: fund name (plain text, never numeric-looking).
    $q4.Cells.Item($r, 3).Value = $data[1]
    $q4.Cells.Item($r, 3).Style = "Normal"

    # Columns D, E, F, G: numeric-looking figures stored as text.
    for ($col = 4; $col -le 7; $col++) {
        $q4.Cells.Item($r, $col).Value = "'" + $data[$col - 2]
        $q4.Cells.Item($r, $col).Style = "Normal"
    }

    # Column H: rank, a genuine number.
    $q4.Cells.Item($r, 8).Value = $data[6]
}
